# removed tool owners pii
# Replace the real usernames in the "Owners" column of the ADS sheet with
# first-name pseudonyms, and re-point the saved "active" selection to the
# ADS sheet (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADS")

$ws.Range("D2").Value  = "Nick"
$ws.Range("D3").Value  = "Milan"
$ws.Range("D4").Value  = "Tolu"
$ws.Range("D5").Value  = "Bousso; Tolu"
$ws.Range("D6").Value  = "Marius"
$ws.Range("D7").Value  = "Milan"
$ws.Range("D8").Value  = "Marius; Milan"
$ws.Range("D9").Value  = "Wayne"
$ws.Range("D10").Value = "Milan"

# Make ADS the active/selected sheet (was NET) with D11 selected, mirroring
# the workbook's saved view state after the edit.
$ws.Select()
$ws.Range("D11").Select()
